$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing 12 rows (134-145) that no longer exist after the edit
$ws.Range("A134:A145").EntireRow.Delete()

# Overwrite rows 2-133 with the refreshed data pull
$ws.Cells.Item(2,1).Value = 45922
$ws.Cells.Item(2,2).Value = 0.021
$ws.Cells.Item(2,3).Value = 8.951
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = "22.09.20251"
$ws.Cells.Item(3,1).Value = 45922.01041666666
$ws.Cells.Item(3,2).Value = 0.006
$ws.Cells.Item(3,3).Value = 1.843
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 2
$ws.Cells.Item(3,7).Value = "22.09.20252"
$ws.Cells.Item(4,1).Value = 45922.02083333334
$ws.Cells.Item(4,2).Value = 0
$ws.Cells.Item(4,3).Value = 4.514
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 3
$ws.Cells.Item(4,7).Value = "22.09.20253"
$ws.Cells.Item(5,1).Value = 45922.03125
$ws.Cells.Item(5,2).Value = 0
$ws.Cells.Item(5,3).Value = 8.714
$ws.Cells.Item(5,4).Value = 0
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 4
$ws.Cells.Item(5,7).Value = "22.09.20254"
$ws.Cells.Item(6,1).Value = 45922.04166666666
$ws.Cells.Item(6,2).Value = 0
$ws.Cells.Item(6,3).Value = 9.493
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = 5
$ws.Cells.Item(6,7).Value = "22.09.20255"
$ws.Cells.Item(7,1).Value = 45922.05208333334
$ws.Cells.Item(7,2).Value = 0
$ws.Cells.Item(7,3).Value = 12.045
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = 0
$ws.Cells.Item(7,6).Value = 6
$ws.Cells.Item(7,7).Value = "22.09.20256"
$ws.Cells.Item(8,1).Value = 45922.0625
$ws.Cells.Item(8,2).Value = 0
$ws.Cells.Item(8,3).Value = 2.642
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = 0
$ws.Cells.Item(8,6).Value = 7
$ws.Cells.Item(8,7).Value = "22.09.20257"
$ws.Cells.Item(9,1).Value = 45922.07291666666
$ws.Cells.Item(9,2).Value = 0
$ws.Cells.Item(9,3).Value = 9.838
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = 0
$ws.Cells.Item(9,6).Value = 8
$ws.Cells.Item(9,7).Value = "22.09.20258"
$ws.Cells.Item(10,1).Value = 45922.08333333334
$ws.Cells.Item(10,2).Value = 0.076
$ws.Cells.Item(10,3).Value = 1.416
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).Value = 0
$ws.Cells.Item(10,6).Value = 9
$ws.Cells.Item(10,7).Value = "22.09.20259"
$ws.Cells.Item(11,1).Value = 45922.09375
$ws.Cells.Item(11,2).Value = 0
$ws.Cells.Item(11,3).Value = 19.579
$ws.Cells.Item(11,4).Value = 0
$ws.Cells.Item(11,5).Value = 0
$ws.Cells.Item(11,6).Value = 10
$ws.Cells.Item(11,7).Value = "22.09.202510"
$ws.Cells.Item(12,1).Value = 45922.10416666666
$ws.Cells.Item(12,2).Value = 0
$ws.Cells.Item(12,3).Value = 7.226
$ws.Cells.Item(12,4).Value = 0
$ws.Cells.Item(12,5).Value = 0
$ws.Cells.Item(12,6).Value = 11
$ws.Cells.Item(12,7).Value = "22.09.202511"
$ws.Cells.Item(13,1).Value = 45922.11458333334
$ws.Cells.Item(13,2).Value = 0
$ws.Cells.Item(13,3).Value = 3.693
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(13,6).Value = 12
$ws.Cells.Item(13,7).Value = "22.09.202512"
$ws.Cells.Item(14,1).Value = 45922.125
$ws.Cells.Item(14,2).Value = 0
$ws.Cells.Item(14,3).Value = 8.037
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 0
$ws.Cells.Item(14,6).Value = 13
$ws.Cells.Item(14,7).Value = "22.09.202513"
$ws.Cells.Item(15,1).Value = 45922.13541666666
$ws.Cells.Item(15,2).Value = 0
$ws.Cells.Item(15,3).Value = 13.317
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 0
$ws.Cells.Item(15,6).Value = 14
$ws.Cells.Item(15,7).Value = "22.09.202514"
$ws.Cells.Item(16,1).Value = 45922.14583333334
$ws.Cells.Item(16,2).Value = 0.241
$ws.Cells.Item(16,3).Value = 3.336
$ws.Cells.Item(16,4).Value = 0
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(16,6).Value = 15
$ws.Cells.Item(16,7).Value = "22.09.202515"
$ws.Cells.Item(17,1).Value = 45922.15625
$ws.Cells.Item(17,2).Value = 0.161
$ws.Cells.Item(17,3).Value = 0.643
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = 0
$ws.Cells.Item(17,6).Value = 16
$ws.Cells.Item(17,7).Value = "22.09.202516"
$ws.Cells.Item(18,1).Value = 45922.16666666666
$ws.Cells.Item(18,2).Value = 0.081
$ws.Cells.Item(18,3).Value = 2.177
$ws.Cells.Item(18,4).Value = 0
$ws.Cells.Item(18,5).Value = 0
$ws.Cells.Item(18,6).Value = 17
$ws.Cells.Item(18,7).Value = "22.09.202517"
$ws.Cells.Item(19,1).Value = 45922.17708333334
$ws.Cells.Item(19,2).Value = 0
$ws.Cells.Item(19,3).Value = 19.315
$ws.Cells.Item(19,4).Value = 0
$ws.Cells.Item(19,5).Value = 0
$ws.Cells.Item(19,6).Value = 18
$ws.Cells.Item(19,7).Value = "22.09.202518"
$ws.Cells.Item(20,1).Value = 45922.1875
$ws.Cells.Item(20,2).Value = 0.2
$ws.Cells.Item(20,3).Value = 1.667
$ws.Cells.Item(20,4).Value = 0
$ws.Cells.Item(20,5).Value = 0
$ws.Cells.Item(20,6).Value = 19
$ws.Cells.Item(20,7).Value = "22.09.202519"
$ws.Cells.Item(21,1).Value = 45922.19791666666
$ws.Cells.Item(21,2).Value = 5.693
$ws.Cells.Item(21,3).Value = 0.001
$ws.Cells.Item(21,4).Value = 0
$ws.Cells.Item(21,5).Value = 0
$ws.Cells.Item(21,6).Value = 20
$ws.Cells.Item(21,7).Value = "22.09.202520"
$ws.Cells.Item(22,1).Value = 45922.20833333334
$ws.Cells.Item(22,2).Value = 0.229
$ws.Cells.Item(22,3).Value = 0.157
$ws.Cells.Item(22,4).Value = 0
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(22,6).Value = 21
$ws.Cells.Item(22,7).Value = "22.09.202521"
$ws.Cells.Item(23,1).Value = 45922.21875
$ws.Cells.Item(23,2).Value = 0.034
$ws.Cells.Item(23,3).Value = 2.611
$ws.Cells.Item(23,4).Value = 0
$ws.Cells.Item(23,5).Value = 0
$ws.Cells.Item(23,6).Value = 22
$ws.Cells.Item(23,7).Value = "22.09.202522"
$ws.Cells.Item(24,1).Value = 45922.22916666666
$ws.Cells.Item(24,2).Value = 13.363
$ws.Cells.Item(24,3).Value = 0.027
$ws.Cells.Item(24,4).Value = 0
$ws.Cells.Item(24,5).Value = 0
$ws.Cells.Item(24,6).Value = 23
$ws.Cells.Item(24,7).Value = "22.09.202523"
$ws.Cells.Item(25,1).Value = 45922.23958333334
$ws.Cells.Item(25,2).Value = 13.281
$ws.Cells.Item(25,3).Value = 0
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(25,6).Value = 24
$ws.Cells.Item(25,7).Value = "22.09.202524"
$ws.Cells.Item(26,1).Value = 45922.25
$ws.Cells.Item(26,2).Value = 10.262
$ws.Cells.Item(26,3).Value = 0.061
$ws.Cells.Item(26,4).Value = 0
$ws.Cells.Item(26,5).Value = 0
$ws.Cells.Item(26,6).Value = 25
$ws.Cells.Item(26,7).Value = "22.09.202525"
$ws.Cells.Item(27,1).Value = 45922.26041666666
$ws.Cells.Item(27,2).Value = 14.615
$ws.Cells.Item(27,3).Value = 0
$ws.Cells.Item(27,4).Value = 25
$ws.Cells.Item(27,5).Value = 0
$ws.Cells.Item(27,6).Value = 26
$ws.Cells.Item(27,7).Value = "22.09.202526"
$ws.Cells.Item(28,1).Value = 45922.27083333334
$ws.Cells.Item(28,2).Value = 16.63
$ws.Cells.Item(28,3).Value = 0
$ws.Cells.Item(28,4).Value = 25
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(28,6).Value = 27
$ws.Cells.Item(28,7).Value = "22.09.202527"
$ws.Cells.Item(29,1).Value = 45922.28125
$ws.Cells.Item(29,2).Value = 3.071
$ws.Cells.Item(29,3).Value = 0.02
$ws.Cells.Item(29,4).Value = 25
$ws.Cells.Item(29,5).Value = 0
$ws.Cells.Item(29,6).Value = 28
$ws.Cells.Item(29,7).Value = "22.09.202528"
$ws.Cells.Item(30,1).Value = 45922.29166666666
$ws.Cells.Item(30,2).Value = 0.172
$ws.Cells.Item(30,3).Value = 0.43
$ws.Cells.Item(30,4).Value = 25
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(30,6).Value = 29
$ws.Cells.Item(30,7).Value = "22.09.202529"
$ws.Cells.Item(31,1).Value = 45922.30208333334
$ws.Cells.Item(31,2).Value = 0
$ws.Cells.Item(31,3).Value = 3.465
$ws.Cells.Item(31,4).Value = 25
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(31,6).Value = 30
$ws.Cells.Item(31,7).Value = "22.09.202530"
$ws.Cells.Item(32,1).Value = 45922.3125
$ws.Cells.Item(32,2).Value = 0
$ws.Cells.Item(32,3).Value = 16.459
$ws.Cells.Item(32,4).Value = 25
$ws.Cells.Item(32,5).Value = 0
$ws.Cells.Item(32,6).Value = 31
$ws.Cells.Item(32,7).Value = "22.09.202531"
$ws.Cells.Item(33,1).Value = 45922.32291666666
$ws.Cells.Item(33,2).Value = 0
$ws.Cells.Item(33,3).Value = 18.164
$ws.Cells.Item(33,4).Value = 0
$ws.Cells.Item(33,5).Value = 0
$ws.Cells.Item(33,6).Value = 32
$ws.Cells.Item(33,7).Value = "22.09.202532"
$ws.Cells.Item(34,1).Value = 45922.33333333334
$ws.Cells.Item(34,2).Value = 10.309
$ws.Cells.Item(34,3).Value = 1.497
$ws.Cells.Item(34,4).Value = 0
$ws.Cells.Item(34,5).Value = 0
$ws.Cells.Item(34,6).Value = 33
$ws.Cells.Item(34,7).Value = "22.09.202533"
$ws.Cells.Item(35,1).Value = 45922.34375
$ws.Cells.Item(35,2).Value = 0.525
$ws.Cells.Item(35,3).Value = 0.519
$ws.Cells.Item(35,4).Value = 0
$ws.Cells.Item(35,5).Value = 0
$ws.Cells.Item(35,6).Value = 34
$ws.Cells.Item(35,7).Value = "22.09.202534"
$ws.Cells.Item(36,1).Value = 45922.35416666666
$ws.Cells.Item(36,2).Value = 0
$ws.Cells.Item(36,3).Value = 32.2
$ws.Cells.Item(36,4).Value = 0
$ws.Cells.Item(36,5).Value = 0
$ws.Cells.Item(36,6).Value = 35
$ws.Cells.Item(36,7).Value = "22.09.202535"
$ws.Cells.Item(37,1).Value = 45922.36458333334
$ws.Cells.Item(37,2).Value = 0
$ws.Cells.Item(37,3).Value = 71.009
$ws.Cells.Item(37,4).Value = 0
$ws.Cells.Item(37,5).Value = 0
$ws.Cells.Item(37,6).Value = 36
$ws.Cells.Item(37,7).Value = "22.09.202536"
$ws.Cells.Item(38,1).Value = 45922.375
$ws.Cells.Item(38,2).Value = 0.314
$ws.Cells.Item(38,3).Value = 7.222
$ws.Cells.Item(38,4).Value = 0
$ws.Cells.Item(38,5).Value = 0
$ws.Cells.Item(38,6).Value = 37
$ws.Cells.Item(38,7).Value = "22.09.202537"
$ws.Cells.Item(39,1).Value = 45922.38541666666
$ws.Cells.Item(39,2).Value = 0
$ws.Cells.Item(39,3).Value = 3.521
$ws.Cells.Item(39,4).Value = 0
$ws.Cells.Item(39,5).Value = 0
$ws.Cells.Item(39,6).Value = 38
$ws.Cells.Item(39,7).Value = "22.09.202538"
$ws.Cells.Item(40,1).Value = 45922.39583333334
$ws.Cells.Item(40,2).Value = 0
$ws.Cells.Item(40,3).Value = 12.614
$ws.Cells.Item(40,4).Value = 0
$ws.Cells.Item(40,5).Value = 0
$ws.Cells.Item(40,6).Value = 39
$ws.Cells.Item(40,7).Value = "22.09.202539"
$ws.Cells.Item(41,1).Value = 45922.40625
$ws.Cells.Item(41,2).Value = 0.492
$ws.Cells.Item(41,3).Value = 4.083
$ws.Cells.Item(41,4).Value = 0
$ws.Cells.Item(41,5).Value = 0
$ws.Cells.Item(41,6).Value = 40
$ws.Cells.Item(41,7).Value = "22.09.202540"
$ws.Cells.Item(42,1).Value = 45922.41666666666
$ws.Cells.Item(42,2).Value = 0.002
$ws.Cells.Item(42,3).Value = 7.928
$ws.Cells.Item(42,4).Value = 0
$ws.Cells.Item(42,5).Value = 0
$ws.Cells.Item(42,6).Value = 41
$ws.Cells.Item(42,7).Value = "22.09.202541"
$ws.Cells.Item(43,1).Value = 45922.42708333334
$ws.Cells.Item(43,2).Value = 0
$ws.Cells.Item(43,3).Value = 2.125
$ws.Cells.Item(43,4).Value = 0
$ws.Cells.Item(43,5).Value = 0
$ws.Cells.Item(43,6).Value = 42
$ws.Cells.Item(43,7).Value = "22.09.202542"
$ws.Cells.Item(44,1).Value = 45922.4375
$ws.Cells.Item(44,2).Value = 0.01
$ws.Cells.Item(44,3).Value = 6.503
$ws.Cells.Item(44,4).Value = 0
$ws.Cells.Item(44,5).Value = 0
$ws.Cells.Item(44,6).Value = 43
$ws.Cells.Item(44,7).Value = "22.09.202543"
$ws.Cells.Item(45,1).Value = 45922.44791666666
$ws.Cells.Item(45,2).Value = 0
$ws.Cells.Item(45,3).Value = 23.516
$ws.Cells.Item(45,4).Value = 0
$ws.Cells.Item(45,5).Value = 0
$ws.Cells.Item(45,6).Value = 44
$ws.Cells.Item(45,7).Value = "22.09.202544"
$ws.Cells.Item(46,1).Value = 45922.45833333334
$ws.Cells.Item(46,2).Value = 0
$ws.Cells.Item(46,3).Value = 21.495
$ws.Cells.Item(46,4).Value = 0
$ws.Cells.Item(46,5).Value = 0
$ws.Cells.Item(46,6).Value = 45
$ws.Cells.Item(46,7).Value = "22.09.202545"
$ws.Cells.Item(47,1).Value = 45922.46875
$ws.Cells.Item(47,2).Value = 0
$ws.Cells.Item(47,3).Value = 9.557
$ws.Cells.Item(47,4).Value = 0
$ws.Cells.Item(47,5).Value = 0
$ws.Cells.Item(47,6).Value = 46
$ws.Cells.Item(47,7).Value = "22.09.202546"
$ws.Cells.Item(48,1).Value = 45922.47916666666
$ws.Cells.Item(48,2).Value = 0
$ws.Cells.Item(48,3).Value = 16.393
$ws.Cells.Item(48,4).Value = 0
$ws.Cells.Item(48,5).Value = 0
$ws.Cells.Item(48,6).Value = 47
$ws.Cells.Item(48,7).Value = "22.09.202547"
$ws.Cells.Item(49,1).Value = 45922.48958333334
$ws.Cells.Item(49,2).Value = 0
$ws.Cells.Item(49,3).Value = 3.601
$ws.Cells.Item(49,4).Value = 0
$ws.Cells.Item(49,5).Value = 0
$ws.Cells.Item(49,6).Value = 48
$ws.Cells.Item(49,7).Value = "22.09.202548"
$ws.Cells.Item(50,1).Value = 45922.5
$ws.Cells.Item(50,2).Value = 0
$ws.Cells.Item(50,3).Value = 12.38
$ws.Cells.Item(50,4).Value = 0
$ws.Cells.Item(50,5).Value = 0
$ws.Cells.Item(50,6).Value = 49
$ws.Cells.Item(50,7).Value = "22.09.202549"
$ws.Cells.Item(51,1).Value = 45922.51041666666
$ws.Cells.Item(51,2).Value = 0
$ws.Cells.Item(51,3).Value = 6.467
$ws.Cells.Item(51,4).Value = 0
$ws.Cells.Item(51,5).Value = 0
$ws.Cells.Item(51,6).Value = 50
$ws.Cells.Item(51,7).Value = "22.09.202550"
$ws.Cells.Item(52,1).Value = 45922.52083333334
$ws.Cells.Item(52,2).Value = 0.002
$ws.Cells.Item(52,3).Value = 6.857
$ws.Cells.Item(52,4).Value = 0
$ws.Cells.Item(52,5).Value = 0
$ws.Cells.Item(52,6).Value = 51
$ws.Cells.Item(52,7).Value = "22.09.202551"
$ws.Cells.Item(53,1).Value = 45922.53125
$ws.Cells.Item(53,2).Value = 0
$ws.Cells.Item(53,3).Value = 7.991
$ws.Cells.Item(53,4).Value = 0
$ws.Cells.Item(53,5).Value = 0
$ws.Cells.Item(53,6).Value = 52
$ws.Cells.Item(53,7).Value = "22.09.202552"
$ws.Cells.Item(54,1).Value = 45922.54166666666
$ws.Cells.Item(54,2).Value = 0
$ws.Cells.Item(54,3).Value = 28.253
$ws.Cells.Item(54,4).Value = 0
$ws.Cells.Item(54,5).Value = 0
$ws.Cells.Item(54,6).Value = 53
$ws.Cells.Item(54,7).Value = "22.09.202553"
$ws.Cells.Item(55,1).Value = 45922.55208333334
$ws.Cells.Item(55,2).Value = 0
$ws.Cells.Item(55,3).Value = 12.413
$ws.Cells.Item(55,4).Value = 0
$ws.Cells.Item(55,5).Value = 0
$ws.Cells.Item(55,6).Value = 54
$ws.Cells.Item(55,7).Value = "22.09.202554"
$ws.Cells.Item(56,1).Value = 45922.5625
$ws.Cells.Item(56,2).Value = 0.03
$ws.Cells.Item(56,3).Value = 2.756
$ws.Cells.Item(56,4).Value = 0
$ws.Cells.Item(56,5).Value = 0
$ws.Cells.Item(56,6).Value = 55
$ws.Cells.Item(56,7).Value = "22.09.202555"
$ws.Cells.Item(57,1).Value = 45922.57291666666
$ws.Cells.Item(57,2).Value = 0.022
$ws.Cells.Item(57,3).Value = 5.237
$ws.Cells.Item(57,4).Value = 0
$ws.Cells.Item(57,5).Value = 0
$ws.Cells.Item(57,6).Value = 56
$ws.Cells.Item(57,7).Value = "22.09.202556"
$ws.Cells.Item(58,1).Value = 45922.58333333334
$ws.Cells.Item(58,2).Value = 0.027
$ws.Cells.Item(58,3).Value = 7.389
$ws.Cells.Item(58,4).Value = 0
$ws.Cells.Item(58,5).Value = 0
$ws.Cells.Item(58,6).Value = 57
$ws.Cells.Item(58,7).Value = "22.09.202557"
$ws.Cells.Item(59,1).Value = 45922.59375
$ws.Cells.Item(59,2).Value = 0.019
$ws.Cells.Item(59,3).Value = 1.241
$ws.Cells.Item(59,4).Value = 0
$ws.Cells.Item(59,5).Value = 0
$ws.Cells.Item(59,6).Value = 58
$ws.Cells.Item(59,7).Value = "22.09.202558"
$ws.Cells.Item(60,1).Value = 45922.60416666666
$ws.Cells.Item(60,2).Value = 0.041
$ws.Cells.Item(60,3).Value = 9.814
$ws.Cells.Item(60,4).Value = 0
$ws.Cells.Item(60,5).Value = 0
$ws.Cells.Item(60,6).Value = 59
$ws.Cells.Item(60,7).Value = "22.09.202559"
$ws.Cells.Item(61,1).Value = 45922.61458333334
$ws.Cells.Item(61,2).Value = 0.983
$ws.Cells.Item(61,3).Value = 0.06
$ws.Cells.Item(61,4).Value = 0
$ws.Cells.Item(61,5).Value = 0
$ws.Cells.Item(61,6).Value = 60
$ws.Cells.Item(61,7).Value = "22.09.202560"
$ws.Cells.Item(62,1).Value = 45922.625
$ws.Cells.Item(62,2).Value = 0
$ws.Cells.Item(62,3).Value = 29.529
$ws.Cells.Item(62,4).Value = 0
$ws.Cells.Item(62,5).Value = 0
$ws.Cells.Item(62,6).Value = 61
$ws.Cells.Item(62,7).Value = "22.09.202561"
$ws.Cells.Item(63,1).Value = 45922.63541666666
$ws.Cells.Item(63,2).Value = 0.005
$ws.Cells.Item(63,3).Value = 8.884
$ws.Cells.Item(63,4).Value = 0
$ws.Cells.Item(63,5).Value = 0
$ws.Cells.Item(63,6).Value = 62
$ws.Cells.Item(63,7).Value = "22.09.202562"
$ws.Cells.Item(64,1).Value = 45922.64583333334
$ws.Cells.Item(64,2).Value = 0.029
$ws.Cells.Item(64,3).Value = 3.002
$ws.Cells.Item(64,4).Value = 0
$ws.Cells.Item(64,5).Value = 0
$ws.Cells.Item(64,6).Value = 63
$ws.Cells.Item(64,7).Value = "22.09.202563"
$ws.Cells.Item(65,1).Value = 45922.65625
$ws.Cells.Item(65,2).Value = 7.005
$ws.Cells.Item(65,3).Value = 0.109
$ws.Cells.Item(65,4).Value = 0
$ws.Cells.Item(65,5).Value = 0
$ws.Cells.Item(65,6).Value = 64
$ws.Cells.Item(65,7).Value = "22.09.202564"
$ws.Cells.Item(66,1).Value = 45922.66666666666
$ws.Cells.Item(66,2).Value = 0
$ws.Cells.Item(66,3).Value = 32.908
$ws.Cells.Item(66,4).Value = 0
$ws.Cells.Item(66,5).Value = 0
$ws.Cells.Item(66,6).Value = 65
$ws.Cells.Item(66,7).Value = "22.09.202565"
$ws.Cells.Item(67,1).Value = 45922.67708333334
$ws.Cells.Item(67,2).Value = 0.011
$ws.Cells.Item(67,3).Value = 10.824
$ws.Cells.Item(67,4).Value = 0
$ws.Cells.Item(67,5).Value = 0
$ws.Cells.Item(67,6).Value = 66
$ws.Cells.Item(67,7).Value = "22.09.202566"
$ws.Cells.Item(68,1).Value = 45922.6875
$ws.Cells.Item(68,2).Value = 3.222
$ws.Cells.Item(68,3).Value = 0.084
$ws.Cells.Item(68,4).Value = 0
$ws.Cells.Item(68,5).Value = 0
$ws.Cells.Item(68,6).Value = 67
$ws.Cells.Item(68,7).Value = "22.09.202567"
$ws.Cells.Item(69,1).Value = 45922.69791666666
$ws.Cells.Item(69,2).Value = 39.461
$ws.Cells.Item(69,3).Value = 0
$ws.Cells.Item(69,4).Value = 0
$ws.Cells.Item(69,5).Value = 0
$ws.Cells.Item(69,6).Value = 68
$ws.Cells.Item(69,7).Value = "22.09.202568"
$ws.Cells.Item(70,1).Value = 45922.70833333334
$ws.Cells.Item(70,2).Value = 2.168
$ws.Cells.Item(70,3).Value = 15.564
$ws.Cells.Item(70,4).Value = 0
$ws.Cells.Item(70,5).Value = 0
$ws.Cells.Item(70,6).Value = 69
$ws.Cells.Item(70,7).Value = "22.09.202569"
$ws.Cells.Item(71,1).Value = 45922.71875
$ws.Cells.Item(71,2).Value = 1.88
$ws.Cells.Item(71,3).Value = 7.269
$ws.Cells.Item(71,4).Value = 0
$ws.Cells.Item(71,5).Value = 0
$ws.Cells.Item(71,6).Value = 70
$ws.Cells.Item(71,7).Value = "22.09.202570"
$ws.Cells.Item(72,1).Value = 45922.72916666666
$ws.Cells.Item(72,2).Value = 11.279
$ws.Cells.Item(72,3).Value = 0
$ws.Cells.Item(72,4).Value = 0
$ws.Cells.Item(72,5).Value = 0
$ws.Cells.Item(72,6).Value = 71
$ws.Cells.Item(72,7).Value = "22.09.202571"
$ws.Cells.Item(73,1).Value = 45922.73958333334
$ws.Cells.Item(73,2).Value = 46.769
$ws.Cells.Item(73,3).Value = 0
$ws.Cells.Item(73,4).Value = 0
$ws.Cells.Item(73,5).Value = 0
$ws.Cells.Item(73,6).Value = 72
$ws.Cells.Item(73,7).Value = "22.09.202572"
$ws.Cells.Item(74,1).Value = 45922.75
$ws.Cells.Item(74,2).Value = 8.482
$ws.Cells.Item(74,3).Value = 0.008
$ws.Cells.Item(74,4).Value = 0
$ws.Cells.Item(74,5).Value = 0
$ws.Cells.Item(74,6).Value = 73
$ws.Cells.Item(74,7).Value = "22.09.202573"
$ws.Cells.Item(75,1).Value = 45922.76041666666
$ws.Cells.Item(75,2).Value = 20.185
$ws.Cells.Item(75,3).Value = 0
$ws.Cells.Item(75,4).Value = 0
$ws.Cells.Item(75,5).Value = 0
$ws.Cells.Item(75,6).Value = 74
$ws.Cells.Item(75,7).Value = "22.09.202574"
$ws.Cells.Item(76,1).Value = 45922.77083333334
$ws.Cells.Item(76,2).Value = 26.5
$ws.Cells.Item(76,3).Value = 0
$ws.Cells.Item(76,4).Value = 0
$ws.Cells.Item(76,5).Value = 0
$ws.Cells.Item(76,6).Value = 75
$ws.Cells.Item(76,7).Value = "22.09.202575"
$ws.Cells.Item(77,1).Value = 45922.78125
$ws.Cells.Item(77,2).Value = 32.723
$ws.Cells.Item(77,3).Value = 0
$ws.Cells.Item(77,4).Value = 25
$ws.Cells.Item(77,5).Value = 0
$ws.Cells.Item(77,6).Value = 76
$ws.Cells.Item(77,7).Value = "22.09.202576"
$ws.Cells.Item(78,1).Value = 45922.79166666666
$ws.Cells.Item(78,2).Value = 4.783
$ws.Cells.Item(78,3).Value = 0.004
$ws.Cells.Item(78,4).Value = 25
$ws.Cells.Item(78,5).Value = 0
$ws.Cells.Item(78,6).Value = 77
$ws.Cells.Item(78,7).Value = "22.09.202577"
$ws.Cells.Item(79,1).Value = 45922.80208333334
$ws.Cells.Item(79,2).Value = 5.662
$ws.Cells.Item(79,3).Value = 0
$ws.Cells.Item(79,4).Value = 0
$ws.Cells.Item(79,5).Value = 0
$ws.Cells.Item(79,6).Value = 78
$ws.Cells.Item(79,7).Value = "22.09.202578"
$ws.Cells.Item(80,1).Value = 45922.8125
$ws.Cells.Item(80,2).Value = 1.831
$ws.Cells.Item(80,3).Value = 0.029
$ws.Cells.Item(80,4).Value = 0
$ws.Cells.Item(80,5).Value = 0
$ws.Cells.Item(80,6).Value = 79
$ws.Cells.Item(80,7).Value = "22.09.202579"
$ws.Cells.Item(81,1).Value = 45922.82291666666
$ws.Cells.Item(81,2).Value = 0.047
$ws.Cells.Item(81,3).Value = 0.317
$ws.Cells.Item(81,4).Value = 0
$ws.Cells.Item(81,5).Value = 0
$ws.Cells.Item(81,6).Value = 80
$ws.Cells.Item(81,7).Value = "22.09.202580"
$ws.Cells.Item(82,1).Value = 45922.83333333334
$ws.Cells.Item(82,2).Value = 1.849
$ws.Cells.Item(82,3).Value = 2.266
$ws.Cells.Item(82,4).Value = 0
$ws.Cells.Item(82,5).Value = 0
$ws.Cells.Item(82,6).Value = 81
$ws.Cells.Item(82,7).Value = "22.09.202581"
$ws.Cells.Item(83,1).Value = 45922.84375
$ws.Cells.Item(83,2).Value = 0
$ws.Cells.Item(83,3).Value = 7.411
$ws.Cells.Item(83,4).Value = 0
$ws.Cells.Item(83,5).Value = 0
$ws.Cells.Item(83,6).Value = 82
$ws.Cells.Item(83,7).Value = "22.09.202582"
$ws.Cells.Item(84,1).Value = 45922.85416666666
$ws.Cells.Item(84,2).Value = 0
$ws.Cells.Item(84,3).Value = 20.421
$ws.Cells.Item(84,4).Value = 0
$ws.Cells.Item(84,5).Value = 0
$ws.Cells.Item(84,6).Value = 83
$ws.Cells.Item(84,7).Value = "22.09.202583"
$ws.Cells.Item(85,1).Value = 45922.86458333334
$ws.Cells.Item(85,2).Value = 0
$ws.Cells.Item(85,3).Value = 12.859
$ws.Cells.Item(85,4).Value = 0
$ws.Cells.Item(85,5).Value = 0
$ws.Cells.Item(85,6).Value = 84
$ws.Cells.Item(85,7).Value = "22.09.202584"
$ws.Cells.Item(86,1).Value = 45922.875
$ws.Cells.Item(86,2).Value = 3.257
$ws.Cells.Item(86,3).Value = 0.933
$ws.Cells.Item(86,4).Value = 0
$ws.Cells.Item(86,5).Value = 0
$ws.Cells.Item(86,6).Value = 85
$ws.Cells.Item(86,7).Value = "22.09.202585"
$ws.Cells.Item(87,1).Value = 45922.88541666666
$ws.Cells.Item(87,2).Value = 0.363
$ws.Cells.Item(87,3).Value = 2.827
$ws.Cells.Item(87,4).Value = 0
$ws.Cells.Item(87,5).Value = 0
$ws.Cells.Item(87,6).Value = 86
$ws.Cells.Item(87,7).Value = "22.09.202586"
$ws.Cells.Item(88,1).Value = 45922.89583333334
$ws.Cells.Item(88,2).Value = 0
$ws.Cells.Item(88,3).Value = 15.015
$ws.Cells.Item(88,4).Value = 0
$ws.Cells.Item(88,5).Value = 0
$ws.Cells.Item(88,6).Value = 87
$ws.Cells.Item(88,7).Value = "22.09.202587"
$ws.Cells.Item(89,1).Value = 45922.90625
$ws.Cells.Item(89,2).Value = 0
$ws.Cells.Item(89,3).Value = 23.511
$ws.Cells.Item(89,4).Value = 0
$ws.Cells.Item(89,5).Value = 0
$ws.Cells.Item(89,6).Value = 88
$ws.Cells.Item(89,7).Value = "22.09.202588"
$ws.Cells.Item(90,1).Value = 45922.91666666666
$ws.Cells.Item(90,2).Value = 1.274
$ws.Cells.Item(90,3).Value = 1.714
$ws.Cells.Item(90,4).Value = 0
$ws.Cells.Item(90,5).Value = 0
$ws.Cells.Item(90,6).Value = 89
$ws.Cells.Item(90,7).Value = "22.09.202589"
$ws.Cells.Item(91,1).Value = 45922.92708333334
$ws.Cells.Item(91,2).Value = 0.036
$ws.Cells.Item(91,3).Value = 1.767
$ws.Cells.Item(91,4).Value = 0
$ws.Cells.Item(91,5).Value = 0
$ws.Cells.Item(91,6).Value = 90
$ws.Cells.Item(91,7).Value = "22.09.202590"
$ws.Cells.Item(92,1).Value = 45922.9375
$ws.Cells.Item(92,2).Value = 0
$ws.Cells.Item(92,3).Value = 8.71
$ws.Cells.Item(92,4).Value = 0
$ws.Cells.Item(92,5).Value = 0
$ws.Cells.Item(92,6).Value = 91
$ws.Cells.Item(92,7).Value = "22.09.202591"
$ws.Cells.Item(93,1).Value = 45922.94791666666
$ws.Cells.Item(93,2).Value = 0
$ws.Cells.Item(93,3).Value = 39.601
$ws.Cells.Item(93,4).Value = 0
$ws.Cells.Item(93,5).Value = 0
$ws.Cells.Item(93,6).Value = 92
$ws.Cells.Item(93,7).Value = "22.09.202592"
$ws.Cells.Item(94,1).Value = 45922.95833333334
$ws.Cells.Item(94,2).Value = 0
$ws.Cells.Item(94,3).Value = 15.793
$ws.Cells.Item(94,4).Value = 0
$ws.Cells.Item(94,5).Value = 0
$ws.Cells.Item(94,6).Value = 93
$ws.Cells.Item(94,7).Value = "22.09.202593"
$ws.Cells.Item(95,1).Value = 45922.96875
$ws.Cells.Item(95,2).Value = 0.005
$ws.Cells.Item(95,3).Value = 3.222
$ws.Cells.Item(95,4).Value = 0
$ws.Cells.Item(95,5).Value = 0
$ws.Cells.Item(95,6).Value = 94
$ws.Cells.Item(95,7).Value = "22.09.202594"
$ws.Cells.Item(96,1).Value = 45922.97916666666
$ws.Cells.Item(96,2).Value = 0.004
$ws.Cells.Item(96,3).Value = 2.352
$ws.Cells.Item(96,4).Value = 0
$ws.Cells.Item(96,5).Value = 0
$ws.Cells.Item(96,6).Value = 95
$ws.Cells.Item(96,7).Value = "22.09.202595"
$ws.Cells.Item(97,1).Value = 45922.98958333334
$ws.Cells.Item(97,2).Value = 0.002
$ws.Cells.Item(97,3).Value = 2.438
$ws.Cells.Item(97,4).Value = 0
$ws.Cells.Item(97,5).Value = 0
$ws.Cells.Item(97,6).Value = 96
$ws.Cells.Item(97,7).Value = "22.09.202596"
$ws.Cells.Item(98,1).Value = 45923
$ws.Cells.Item(98,2).Value = 6.244
$ws.Cells.Item(98,3).Value = 0.29
$ws.Cells.Item(98,4).Value = 0
$ws.Cells.Item(98,5).Value = 0
$ws.Cells.Item(98,6).Value = 1
$ws.Cells.Item(98,7).Value = "23.09.20251"
$ws.Cells.Item(99,1).Value = 45923
$ws.Cells.Item(99,2).Value = 6.244
$ws.Cells.Item(99,3).Value = 0.29
$ws.Cells.Item(99,4).Value = 0
$ws.Cells.Item(99,5).Value = 0
$ws.Cells.Item(99,6).Value = 1
$ws.Cells.Item(99,7).Value = "23.09.20251"
$ws.Cells.Item(100,1).Value = 45923.01041666666
$ws.Cells.Item(100,2).Value = 0.498
$ws.Cells.Item(100,3).Value = 0.031
$ws.Cells.Item(100,4).Value = 0
$ws.Cells.Item(100,5).Value = 0
$ws.Cells.Item(100,6).Value = 2
$ws.Cells.Item(100,7).Value = "23.09.20252"
$ws.Cells.Item(101,1).Value = 45923.01041666666
$ws.Cells.Item(101,2).Value = 0.498
$ws.Cells.Item(101,3).Value = 0.031
$ws.Cells.Item(101,4).Value = 0
$ws.Cells.Item(101,5).Value = 0
$ws.Cells.Item(101,6).Value = 2
$ws.Cells.Item(101,7).Value = "23.09.20252"
$ws.Cells.Item(102,1).Value = 45923.02083333334
$ws.Cells.Item(102,2).Value = 0.123
$ws.Cells.Item(102,3).Value = 0.161
$ws.Cells.Item(102,4).Value = 0
$ws.Cells.Item(102,5).Value = 0
$ws.Cells.Item(102,6).Value = 3
$ws.Cells.Item(102,7).Value = "23.09.20253"
$ws.Cells.Item(103,1).Value = 45923.02083333334
$ws.Cells.Item(103,2).Value = 0.123
$ws.Cells.Item(103,3).Value = 0.161
$ws.Cells.Item(103,4).Value = 0
$ws.Cells.Item(103,5).Value = 0
$ws.Cells.Item(103,6).Value = 3
$ws.Cells.Item(103,7).Value = "23.09.20253"
$ws.Cells.Item(104,1).Value = 45923.03125
$ws.Cells.Item(104,2).Value = 0.101
$ws.Cells.Item(104,3).Value = 0.123
$ws.Cells.Item(104,4).Value = 0
$ws.Cells.Item(104,5).Value = 0
$ws.Cells.Item(104,6).Value = 4
$ws.Cells.Item(104,7).Value = "23.09.20254"
$ws.Cells.Item(105,1).Value = 45923.03125
$ws.Cells.Item(105,2).Value = 0.101
$ws.Cells.Item(105,3).Value = 0.123
$ws.Cells.Item(105,4).Value = 0
$ws.Cells.Item(105,5).Value = 0
$ws.Cells.Item(105,6).Value = 4
$ws.Cells.Item(105,7).Value = "23.09.20254"
$ws.Cells.Item(106,1).Value = 45923.04166666666
$ws.Cells.Item(106,2).Value = 0
$ws.Cells.Item(106,3).Value = 2.846
$ws.Cells.Item(106,4).Value = 0
$ws.Cells.Item(106,5).Value = 0
$ws.Cells.Item(106,6).Value = 5
$ws.Cells.Item(106,7).Value = "23.09.20255"
$ws.Cells.Item(107,1).Value = 45923.04166666666
$ws.Cells.Item(107,2).Value = 0
$ws.Cells.Item(107,3).Value = 2.846
$ws.Cells.Item(107,4).Value = 0
$ws.Cells.Item(107,5).Value = 0
$ws.Cells.Item(107,6).Value = 5
$ws.Cells.Item(107,7).Value = "23.09.20255"
$ws.Cells.Item(108,1).Value = 45923.05208333334
$ws.Cells.Item(108,2).Value = 0.074
$ws.Cells.Item(108,3).Value = 1.735
$ws.Cells.Item(108,4).Value = 0
$ws.Cells.Item(108,5).Value = 0
$ws.Cells.Item(108,6).Value = 6
$ws.Cells.Item(108,7).Value = "23.09.20256"
$ws.Cells.Item(109,1).Value = 45923.05208333334
$ws.Cells.Item(109,2).Value = 0.074
$ws.Cells.Item(109,3).Value = 1.735
$ws.Cells.Item(109,4).Value = 0
$ws.Cells.Item(109,5).Value = 0
$ws.Cells.Item(109,6).Value = 6
$ws.Cells.Item(109,7).Value = "23.09.20256"
$ws.Cells.Item(110,1).Value = 45923.0625
$ws.Cells.Item(110,2).Value = 0.357
$ws.Cells.Item(110,3).Value = 0.2
$ws.Cells.Item(110,4).Value = 0
$ws.Cells.Item(110,5).Value = 0
$ws.Cells.Item(110,6).Value = 7
$ws.Cells.Item(110,7).Value = "23.09.20257"
$ws.Cells.Item(111,1).Value = 45923.0625
$ws.Cells.Item(111,2).Value = 0.357
$ws.Cells.Item(111,3).Value = 0.2
$ws.Cells.Item(111,4).Value = 0
$ws.Cells.Item(111,5).Value = 0
$ws.Cells.Item(111,6).Value = 7
$ws.Cells.Item(111,7).Value = "23.09.20257"
$ws.Cells.Item(112,1).Value = 45923.07291666666
$ws.Cells.Item(112,2).Value = 0.906
$ws.Cells.Item(112,3).Value = 0.257
$ws.Cells.Item(112,4).Value = 0
$ws.Cells.Item(112,5).Value = 0
$ws.Cells.Item(112,6).Value = 8
$ws.Cells.Item(112,7).Value = "23.09.20258"
$ws.Cells.Item(113,1).Value = 45923.07291666666
$ws.Cells.Item(113,2).Value = 0.906
$ws.Cells.Item(113,3).Value = 0.257
$ws.Cells.Item(113,4).Value = 0
$ws.Cells.Item(113,5).Value = 0
$ws.Cells.Item(113,6).Value = 8
$ws.Cells.Item(113,7).Value = "23.09.20258"
$ws.Cells.Item(114,1).Value = 45923.08333333334
$ws.Cells.Item(114,2).Value = 0
$ws.Cells.Item(114,3).Value = 2.423
$ws.Cells.Item(114,4).Value = 0
$ws.Cells.Item(114,5).Value = 0
$ws.Cells.Item(114,6).Value = 9
$ws.Cells.Item(114,7).Value = "23.09.20259"
$ws.Cells.Item(115,1).Value = 45923.09375
$ws.Cells.Item(115,2).Value = 0
$ws.Cells.Item(115,3).Value = 9.889
$ws.Cells.Item(115,4).Value = 0
$ws.Cells.Item(115,5).Value = 0
$ws.Cells.Item(115,6).Value = 10
$ws.Cells.Item(115,7).Value = "23.09.202510"
$ws.Cells.Item(116,1).Value = 45923.10416666666
$ws.Cells.Item(116,2).Value = 0
$ws.Cells.Item(116,3).Value = 4.696
$ws.Cells.Item(116,4).Value = 0
$ws.Cells.Item(116,5).Value = 0
$ws.Cells.Item(116,6).Value = 11
$ws.Cells.Item(116,7).Value = "23.09.202511"
$ws.Cells.Item(117,1).Value = 45923.11458333334
$ws.Cells.Item(117,2).Value = 0.025
$ws.Cells.Item(117,3).Value = 0.587
$ws.Cells.Item(117,4).Value = 0
$ws.Cells.Item(117,5).Value = 0
$ws.Cells.Item(117,6).Value = 12
$ws.Cells.Item(117,7).Value = "23.09.202512"
$ws.Cells.Item(118,1).Value = 45923.125
$ws.Cells.Item(118,2).Value = 2.476
$ws.Cells.Item(118,3).Value = 0.04
$ws.Cells.Item(118,4).Value = 0
$ws.Cells.Item(118,5).Value = 0
$ws.Cells.Item(118,6).Value = 13
$ws.Cells.Item(118,7).Value = "23.09.202513"
$ws.Cells.Item(119,1).Value = 45923.13541666666
$ws.Cells.Item(119,2).Value = 4.203
$ws.Cells.Item(119,3).Value = 0
$ws.Cells.Item(119,4).Value = 0
$ws.Cells.Item(119,5).Value = 0
$ws.Cells.Item(119,6).Value = 14
$ws.Cells.Item(119,7).Value = "23.09.202514"
$ws.Cells.Item(120,1).Value = 45923.14583333334
$ws.Cells.Item(120,2).Value = 0.721
$ws.Cells.Item(120,3).Value = 0.02
$ws.Cells.Item(120,4).Value = 0
$ws.Cells.Item(120,5).Value = 0
$ws.Cells.Item(120,6).Value = 15
$ws.Cells.Item(120,7).Value = "23.09.202515"
$ws.Cells.Item(121,1).Value = 45923.15625
$ws.Cells.Item(121,2).Value = 0.048
$ws.Cells.Item(121,3).Value = 0.875
$ws.Cells.Item(121,4).Value = 0
$ws.Cells.Item(121,5).Value = 0
$ws.Cells.Item(121,6).Value = 16
$ws.Cells.Item(121,7).Value = "23.09.202516"
$ws.Cells.Item(122,1).Value = 45923.16666666666
$ws.Cells.Item(122,2).Value = 0.038
$ws.Cells.Item(122,3).Value = 0.628
$ws.Cells.Item(122,4).Value = 0
$ws.Cells.Item(122,5).Value = 0
$ws.Cells.Item(122,6).Value = 17
$ws.Cells.Item(122,7).Value = "23.09.202517"
$ws.Cells.Item(123,1).Value = 45923.17708333334
$ws.Cells.Item(123,2).Value = 0.002
$ws.Cells.Item(123,3).Value = 2.466
$ws.Cells.Item(123,4).Value = 0
$ws.Cells.Item(123,5).Value = 0
$ws.Cells.Item(123,6).Value = 18
$ws.Cells.Item(123,7).Value = "23.09.202518"
$ws.Cells.Item(124,1).Value = 45923.1875
$ws.Cells.Item(124,2).Value = 0.165
$ws.Cells.Item(124,3).Value = 0.053
$ws.Cells.Item(124,4).Value = 0
$ws.Cells.Item(124,5).Value = 0
$ws.Cells.Item(124,6).Value = 19
$ws.Cells.Item(124,7).Value = "23.09.202519"
$ws.Cells.Item(125,1).Value = 45923.19791666666
$ws.Cells.Item(125,2).Value = 2.426
$ws.Cells.Item(125,3).Value = 0.103
$ws.Cells.Item(125,4).Value = 0
$ws.Cells.Item(125,5).Value = 0
$ws.Cells.Item(125,6).Value = 20
$ws.Cells.Item(125,7).Value = "23.09.202520"
$ws.Cells.Item(126,1).Value = 45923.20833333334
$ws.Cells.Item(126,2).Value = 4.545
$ws.Cells.Item(126,3).Value = 0
$ws.Cells.Item(126,4).Value = 0
$ws.Cells.Item(126,5).Value = 0
$ws.Cells.Item(126,6).Value = 21
$ws.Cells.Item(126,7).Value = "23.09.202521"
$ws.Cells.Item(127,1).Value = 45923.21875
$ws.Cells.Item(127,2).Value = 7.634
$ws.Cells.Item(127,3).Value = 0
$ws.Cells.Item(127,4).Value = 0
$ws.Cells.Item(127,5).Value = 0
$ws.Cells.Item(127,6).Value = 22
$ws.Cells.Item(127,7).Value = "23.09.202522"
$ws.Cells.Item(128,1).Value = 45923.22916666666
$ws.Cells.Item(128,2).Value = 15.913
$ws.Cells.Item(128,3).Value = 0
$ws.Cells.Item(128,4).Value = 0
$ws.Cells.Item(128,5).Value = 0
$ws.Cells.Item(128,6).Value = 23
$ws.Cells.Item(128,7).Value = "23.09.202523"
$ws.Cells.Item(129,1).Value = 45923.23958333334
$ws.Cells.Item(129,2).Value = 19.479
$ws.Cells.Item(129,3).Value = 0
$ws.Cells.Item(129,4).Value = 0
$ws.Cells.Item(129,5).Value = 0
$ws.Cells.Item(129,6).Value = 24
$ws.Cells.Item(129,7).Value = "23.09.202524"
$ws.Cells.Item(130,1).Value = 45923.25
$ws.Cells.Item(130,2).Value = 8.147
$ws.Cells.Item(130,3).Value = 2.996
$ws.Cells.Item(130,4).Value = 0
$ws.Cells.Item(130,5).Value = 0
$ws.Cells.Item(130,6).Value = 25
$ws.Cells.Item(130,7).Value = "23.09.202525"
$ws.Cells.Item(131,1).Value = 45923.26041666666
$ws.Cells.Item(131,2).Value = 26.061
$ws.Cells.Item(131,3).Value = 0
$ws.Cells.Item(131,4).Value = 0
$ws.Cells.Item(131,5).Value = 0
$ws.Cells.Item(131,6).Value = 26
$ws.Cells.Item(131,7).Value = "23.09.202526"
$ws.Cells.Item(132,1).Value = 45923.27083333334
$ws.Cells.Item(132,2).Value = 23.033
$ws.Cells.Item(132,3).Value = 0
$ws.Cells.Item(132,4).Value = 0
$ws.Cells.Item(132,5).Value = 0
$ws.Cells.Item(132,6).Value = 27
$ws.Cells.Item(132,7).Value = "23.09.202527"
$ws.Cells.Item(133,1).Value = 45923.28125
$ws.Cells.Item(133,2).Value = 3.457
$ws.Cells.Item(133,3).Value = 0.002
$ws.Cells.Item(133,4).Value = 0
$ws.Cells.Item(133,5).Value = 0
$ws.Cells.Item(133,6).Value = 28
$ws.Cells.Item(133,7).Value = "23.09.202528"
